$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Materials")

# Remove the "subgenus" column (header in row 1, "${subgenus}" value in row 2)
$ws.Range("AS1:AS2").EntireColumn.Delete()
